$wb = $excel.ActiveWorkbook

# "Custom Types" sheet gets a new row 2: GIFComment maps to the corrected
# xref link used in the documentation table.
$ws = $wb.Worksheets.Item("Custom Types")
$ws.Range("A2").Value = "GIFComment"
$ws.Range("B2").Value = "[GIFComment](xref:ExifLibrary.GIFComment)"
$ws.Range("C2").Value = "string"

# Re-touch the dependent formula chain on "All Tags" so the cached results
# pick up the newly added Custom Types row (N2/P2/Q2/R2 all derive from a
# VLOOKUP against 'Custom Types').
$tags = $wb.Worksheets.Item("All Tags")
$tags.Range("N2").Formula = $tags.Range("N2").Formula

# Reset the stale R2 selection left over on "All Tags" back to A1.
$tags.Activate()
$tags.Range("A1").Select()
